$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing hyperlink cell style (column F data cells use the
# "Hyperlink" style) so it can be re-applied after the hyperlinks are
# rebuilt below.
$linkStyle = $ws.Range("F2").Style

# Remove every existing hyperlink on the sheet. Rows 6-10 (and their
# hyperlinks) are being dropped entirely, and the hyperlinks that remain
# (F2:F5) are going to be re-created with fresh target URLs, so it is
# simplest to clear them all up front.
$ws.Hyperlinks.Delete()

# Drop the old rows 6-10 - only 4 data rows remain after this edit.
$ws.Rows("6:10").Delete()

# Row 2
$ws.Range("A2").Value = "2025-10-13 06:30:06"
$ws.Range("B2").Value = "【急募】クリニック向け内視鏡画像システム開発の依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5412233"
$ws.Range("G2").Value = 125
$ws.Range("H2").Value = "◆開発,システム開発"

# Row 3
$ws.Range("A3").Value = "2025-10-13 06:30:06"
$ws.Range("B3").Value = "【急募】onedrive上のexcelで自動化システム構築依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5412194"
$ws.Range("G3").Value = 95
$ws.Range("H3").Value = "◆自動化"

# Row 4
$ws.Range("A4").Value = "2025-10-13 06:30:06"
$ws.Range("B4").Value = "【急募】スタートアップ向けプロダクト開発のパートナー募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5412179"
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = "◆開発"

# Row 5
$ws.Range("A5").Value = "2025-10-13 06:30:06"
$ws.Range("B5").Value = "微生物の特定と分類を行いたく、画像解析の専門家を探しています!(急いでません!)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5411887"
$ws.Range("G5").Value = 18
$ws.Range("H5").Value = ""

# Re-create the hyperlinks for the 4 remaining rows, pointing at the new
# URLs, then restore the original "Hyperlink" style (Hyperlinks.Add swaps
# in its own near-duplicate style otherwise).
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5412233")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5412194")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5412179")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5411887")
$ws.Range("F2:F5").Style = $linkStyle

# Column width tweaks: column B 57 -> 42 characters, column H 17 -> 12
# characters. This engine's ColumnWidth setter adds ~5/6 of a character
# of padding relative to the stored OOXML "width", so back that out here
# to land on the exact target width.
$ws.Columns.Item(2).ColumnWidth = 42 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6
